$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.Value = "'" + $newValue
    $rng.Style = $origStyle
}

Set-TextValue "D2" "36.854.50"
Set-TextValue "E2" "  +0.35%  "
Set-TextValue "D3" "1.979.46"
Set-TextValue "E3" "  +0.86%  "
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "245.35"
Set-TextValue "E5" "  +0.20%  "
Set-TextValue "D6" "0.632"
Set-TextValue "E6" "  +1.94%  "
Set-TextValue "D7" "61.26"
Set-TextValue "E7" "  +3.78%  "
Set-TextValue "E8" "  -0.02%  "
Set-TextValue "D9" "0.382"
Set-TextValue "E9" "  +2.20%  "
Set-TextValue "E10" "  -1.38%  "
Set-TextValue "E11" "  +0.81%  "
Set-TextValue "D12" "14.54"
Set-TextValue "E12" "  +5.76%  "
Set-TextValue "D13" "0.848"
Set-TextValue "D14" "22.12"
Set-TextValue "E14" "  -1.41%  "
Set-TextValue "D15" "2.268.20"
Set-TextValue "E15" "  +0.84%  "
Set-TextValue "D16" "5.42"
Set-TextValue "E16" "  +2.66%  "
Set-TextValue "D17" "1.972.39"
Set-TextValue "E17" "  +0.22%  "
Set-TextValue "D18" "36.751.03"
Set-TextValue "E18" "  +0.36%  "
Set-TextValue "D19" "70.18"
Set-TextValue "E19" "  +0.58%  "
Set-TextValue "D20" "0.0₃0861"
Set-TextValue "E20" "  -0.12%  "
Set-TextValue "D21" "5.15"
Set-TextValue "E21" "  +1.08%  "
Set-TextValue "D22" "230.48"
Set-TextValue "E22" "  +0.55%  "
Set-TextValue "E23" "  +0.12%  "
Set-TextValue "D24" "2.48"
Set-TextValue "E24" "  +1.64%  "
Set-TextValue "E25" "  +1.07%  "
Set-TextValue "E26" "  +3.48%  "
Set-TextValue "D27" "9.28"
Set-TextValue "E27" "  -0.52%  "
Set-TextValue "D28" "163.17"
Set-TextValue "E28" "  +1.55%  "
Set-TextValue "D29" "19.49"
Set-TextValue "D30" "1.35"
Set-TextValue "E30" "  +20.03%  "
Set-TextValue "E31" "  +2.06%  "
Set-TextValue "D32" "4.84"
Set-TextValue "E32" "  +2.29%  "
Set-TextValue "D33" "0.0622"
Set-TextValue "E33" "  +0.50%  "
Set-TextValue "D34" "4.53"
Set-TextValue "E34" "  +6.23%  "
Set-TextValue "D35" "2.27"
Set-TextValue "E35" "  +0.74%  "
Set-TextValue "E36" "  -0.14%  "
Set-TextValue "D37" "3.36"
Set-TextValue "E37" "  -1.33%  "
Set-TextValue "E38" "  +0.09%  "
Set-TextValue "D39" "5.49"
Set-TextValue "E39" "  -9.81%  "
Set-TextValue "D40" "0.0974"
Set-TextValue "E40" "  -3.01%  "
Set-TextValue "E42" "  +0.79%  "
Set-TextValue "E43" "  -0.26%  "
Set-TextValue "D44" "16.11"
Set-TextValue "D45" "1.372.18"
Set-TextValue "E45" "  +1.01%  "
Set-TextValue "D46" "89.74"
Set-TextValue "E46" "  +2.23%  "
Set-TextValue "D47" "1.04"
Set-TextValue "E47" "  -0.29%  "
Set-TextValue "D48" "7.24"
Set-TextValue "E48" "  +1.14%  "
Set-TextValue "E49" "  -0.53%  "
Set-TextValue "D50" "46.36"
Set-TextValue "E50" "  +6.03%  "
Set-TextValue "D51" "2.160.85"
Set-TextValue "E51" "  +0.94%  "
